$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.709842930011267
$ws.Range("B3").Value = 0.0840523502246763
$ws.Range("B4").Value = 1.084062404365242
$ws.Range("B5").Value = 0.6827811467223019
$ws.Range("B6").Value = 1.461752845699677
$ws.Range("B7").Value = 0.2174853294999878
